$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOB1240")

# Nome / Name
$ws.Range("B4").Value = "Environmental Geotechnics"
$ws.Range("C4").Value = "Environmental Geotechnics"

# Ativação (force text so the date-like string isn't converted to a date serial)
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").Value = "01/01/2022"

# Objetivos
$ws.Range("B10").Value = "Propiciar o entendimento sobre os conceitos básicos, teoria, metodologias de análise e práticas de engenharia em geotecnia ambiental."
$ws.Range("C10").Value = "Propiciar o entendimento sobre os conceitos básicos, teoria, metodologias de análise e práticas de engenharia em geotecnia ambiental."

# Objectives
$ws.Range("B11").Value = "Provide knowledge of the basic concepts, theory,  analysis methodologies and engineering practices in environmental geotechnics.."
$ws.Range("C11").Value = "Provide knowledge of the basic concepts, theory,  analysis methodologies and engineering practices in environmental geotechnics.."

# Programa resumido
$ws.Range("B14").Value = "Principais fatores condicionantes atuantes em processos como erosão, movimento gravitacional de massa e contaminação de materiais geológicos e água subterrânea. Técnicas de investigação, monitoramento, contenção e recuperação em geotecnia ambiental."
$ws.Range("C14").Value = "Principais fatores condicionantes atuantes em processos como erosão, movimento gravitacional de massa e contaminação de materiais geológicos e água subterrânea. Técnicas de investigação, monitoramento, contenção e recuperação em geotecnia ambiental."

# Short syllabus
$ws.Range("B15").Value = "Main aspects acting in processes such as erosion, gravitational mass movements and contamination of geological materials and groundwater. Investigation, monitoring, containment and recuperation techniques in environmental geotechnics."
$ws.Range("C15").Value = "Main aspects acting in processes such as erosion, gravitational mass movements and contamination of geological materials and groundwater. Investigation, monitoring, containment and recuperation techniques in environmental geotechnics."

# Programa
$ws.Range("B16").Value = "Aspectos conceituais e abrangência do tema. Investigação e monitoramento geoambiental. Mapas geoambientais e geotécnicos. Movimentos gravitacionais de massa: conceitos, tipos e mecanismos. Análise de estabilidade de taludes e técnicas de estabilização de encostas naturais. Processos erosivos. Estimativas de perda de solo e técnicas de mitigação de processos erosivos antrópicos. Contaminação de materiais geológicos e água subterrânea. Transporte de contaminantes na zona saturada e não saturada do solo. Barragens de rejeitos e aterros de resíduos. Geossintéticos"
$ws.Range("C16").Value = "Aspectos conceituais e abrangência do tema. Investigação e monitoramento geoambiental. Mapas geoambientais e geotécnicos. Movimentos gravitacionais de massa: conceitos, tipos e mecanismos. Análise de estabilidade de taludes e técnicas de estabilização de encostas naturais. Processos erosivos. Estimativas de perda de solo e técnicas de mitigação de processos erosivos antrópicos. Contaminação de materiais geológicos e água subterrânea. Transporte de contaminantes na zona saturada e não saturada do solo. Barragens de rejeitos e aterros de resíduos. Geossintéticos"

# Syllabus
$ws.Range("B17").Value = "Conceptual aspects and scope of the theme. Geoenvironmental investigation and monitoring. Geoenvironmental and geotechnical maps. Gravitational mass movements: concepts, types and mechanics. Slope stability analysis and natural slope stabilization techniques. Erosive processes. Estimates of soil loss and mitigation techniques for anthropic erosive processes. Contamination of geological materials and groundwater. Transport of contaminants in saturated and unsaturated zone of the soil. Tailings dams and waste landfills. Geosynthetics."
$ws.Range("C17").Value = "Conceptual aspects and scope of the theme. Geoenvironmental investigation and monitoring. Geoenvironmental and geotechnical maps. Gravitational mass movements: concepts, types and mechanics. Slope stability analysis and natural slope stabilization techniques. Erosive processes. Estimates of soil loss and mitigation techniques for anthropic erosive processes. Contamination of geological materials and groundwater. Transport of contaminants in saturated and unsaturated zone of the soil. Tailings dams and waste landfills. Geosynthetics."

# Método
$ws.Range("B19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$ws.Range("C19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."

# Critério
$ws.Range("B20").Value = "Média ponderada de provas  e atividades."
$ws.Range("C20").Value = "Média ponderada de provas  e atividades."

# Norma de recuperação
$ws.Range("B21").Value = "1 (uma) prova escrita"
$ws.Range("C21").Value = "1 (uma) prova escrita"

# Bibliografia
$ws.Range("B22").Value = "Bibliografia básica:CHIOSSI, n. Geologia de Engenharia. Oficina de Textos: São Paulo, 1979. 427p.Bibliografia complementarBOSCOV, M. E. Geotecnia ambiental. Oficina de Textos: São Paulo, 2008. 248p.ZUQUETTE, L. V. Geotecnia ambiental. Elsevier: Rio de Janeiro, 2015. 432p.OLIVEIRA, A. M. S.; JERÔNIMO, J. Geologia de Engenharia e Ambiental, ABGE, 2018. 912 p.BARROW, C. J. Land degradation. Cambridge University Press: New York, 1991. 316p.KHALED, S. B. D. Fundamentos de Engenharia Geotécnica. Cengage Learning. 2020, 630p.SHARMA, H. D.; REDDY, K. R. Geoenvironmental engineering, Wiley, 2004. 992p. YONG, R. N. Geoenvironmental engineering. contaminated soils, Pollutant fate and migration. CRC Press, 2001. 307p.DANIEL, D.E.  Geotechnical practice for waste disposal. Chapman & Hall, 1993. 693p.ROWE, R. K.; QUIGLEY, R.M.; BOOKER, J. R. Clayey barrier systems for waste disposal facilities. E & FN SPON, 2019. 600p."
$ws.Range("C22").Value = "Bibliografia básica:CHIOSSI, n. Geologia de Engenharia. Oficina de Textos: São Paulo, 1979. 427p.Bibliografia complementarBOSCOV, M. E. Geotecnia ambiental. Oficina de Textos: São Paulo, 2008. 248p.ZUQUETTE, L. V. Geotecnia ambiental. Elsevier: Rio de Janeiro, 2015. 432p.OLIVEIRA, A. M. S.; JERÔNIMO, J. Geologia de Engenharia e Ambiental, ABGE, 2018. 912 p.BARROW, C. J. Land degradation. Cambridge University Press: New York, 1991. 316p.KHALED, S. B. D. Fundamentos de Engenharia Geotécnica. Cengage Learning. 2020, 630p.SHARMA, H. D.; REDDY, K. R. Geoenvironmental engineering, Wiley, 2004. 992p. YONG, R. N. Geoenvironmental engineering. contaminated soils, Pollutant fate and migration. CRC Press, 2001. 307p.DANIEL, D.E.  Geotechnical practice for waste disposal. Chapman & Hall, 1993. 693p.ROWE, R. K.; QUIGLEY, R.M.; BOOKER, J. R. Clayey barrier systems for waste disposal facilities. E & FN SPON, 2019. 600p."
